$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F: header + values for the new week (28_01_2024)
$ws.Range("F1").Value = "28_01_2024"
$ws.Range("F2").Value = 2824
$ws.Range("F3").Value = 2097
$ws.Range("F4").Value = 3537
$ws.Range("F5").Value = 6858

# Update selection to reflect the new last empty cell below column F
$ws.Range("F6").Select()
